$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.443.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.590.11"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.93"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.09"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.591.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.531"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.455.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0709"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.411.90"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.61%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.77"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.982"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  +4.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.61"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.730.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0108"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0523"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +16.20%  "
